# Generate Report for Handoff
# - Overview sheet: zh-cn/de-de status goes from "Handed back: in sync with en-US"
#   to "Ready for handoff", and the "Latest HO Xliff Generate Date" timestamp is
#   refreshed.
# - zh-cn / de-de detail sheets: the "Status" cell gets the same refreshed text,
#   and the zh-cn sheet's "Latest Handoff Datetime" timestamp is refreshed too.
# - The Status columns are narrower now that the status text is shorter, so we
#   shrink them to match.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-18 13:00:02"
$wsOverview.Columns("E:F").ColumnWidth = 16.33

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-18 12:59:51"
$wsZhCn.Columns("C:C").ColumnWidth = 16.33

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Columns("C:C").ColumnWidth = 16.33
